# "Docs image with resizing"
# Re-layout the architecture diagram: grow the two account rectangles,
# re-flow/resize the inner icons, labels, connectors and caption boxes,
# and restyle the text runs (new typefaces / sizes, split "Soda cloud
# account" into two differently-styled runs).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# Shape 1 (id 4) "Rectangle 3" - "Soda cloud account"
# ---------------------------------------------------------------------
$shp = $s.Shapes.Item(1)
$shp.Left = 217.82276916503906
$shp.Top = 143.25
$shp.Width = 205.4362335205078
$shp.Height = 225.8639373779297
$shp.TextFrame2.MarginTop = 22.677165354330707

$tr = $shp.TextFrame.TextRange
$tr.Text = "Soda cloud account"
$run1 = $tr.Characters(1, 4)
$run1.Font.Size = 20
$run1.Font.Name = "Suez One"
$run1.Font.NameComplexScript = "Suez One"
$run2 = $tr.Characters(5, 14)
$run2.Font.Size = 20
$run2.Font.Name = "IBM Plex Sans"

# ---------------------------------------------------------------------
# Shape 2 (id 5) "Rectangle 4" - "Your cloud account"
# ---------------------------------------------------------------------
$shp = $s.Shapes.Item(2)
$shp.Left = 514.0253295898438
$shp.Top = 143.25
$shp.Width = 205.063232421875
$shp.Height = 225.8639373779297
$shp.TextFrame2.MarginTop = 22.677165354330707

$tr = $shp.TextFrame.TextRange
$tr.Font.Size = 20
$tr.Font.Name = "IBM Plex Sans"

# ---------------------------------------------------------------------
# Shape 3 (id 6) "Graphic 5" picture
# ---------------------------------------------------------------------
$shp = $s.Shapes.Item(3)
$shp.Left = 535.8291015625
$shp.Top = 212.78843688964844

# ---------------------------------------------------------------------
# Shape 4 (id 7) "Graphic 6" picture
# ---------------------------------------------------------------------
$shp = $s.Shapes.Item(4)
$shp.Left = 535.8291015625
$shp.Top = 289.417724609375

# ---------------------------------------------------------------------
# Shape 5 (id 8) "TextBox 7" - "Your data lake"
# ---------------------------------------------------------------------
$shp = $s.Shapes.Item(5)
$tr = $shp.TextFrame.TextRange
$tr.Font.Size = 14
$tr.Font.Name = "IBM Plex Sans"
$shp.Left = 591.8291015625
$shp.Top = 227.45953369140625
$shp.Width = 105.4195327758789
$shp.Height = 24.234411239624023

# ---------------------------------------------------------------------
# Shape 6 (id 9) "TextBox 8" - "Soda file storage"
# ---------------------------------------------------------------------
$shp = $s.Shapes.Item(6)
$tr = $shp.TextFrame.TextRange
$tr.Font.Size = 14
$tr.Font.Name = "IBM Plex Sans"
$shp.Left = 591.8291015625
$shp.Top = 304.0888366699219
$shp.Width = 120.31354522705078
$shp.Height = 24.234411239624023

# ---------------------------------------------------------------------
# Shape 7 (id 11) "Straight Arrow Connector 10"
# ---------------------------------------------------------------------
$shp = $s.Shapes.Item(7)
$shp.Left = 395.5443420410156
$shp.Top = 241.1862335205078

# ---------------------------------------------------------------------
# Shape 8 (id 16) "TextBox 15" - "read only"
# ---------------------------------------------------------------------
$shp = $s.Shapes.Item(8)
$tr = $shp.TextFrame.TextRange
$tr.Font.Size = 12
$tr.Font.Name = "IBM Plex Sans"
$shp.Left = 437.1703186035156
$shp.Top = 218.20103454589844
$shp.Width = 65.15519714355469
$shp.Height = 21.810945510864258

# ---------------------------------------------------------------------
# Shape 9 (id 19) "Straight Arrow Connector 18"
# ---------------------------------------------------------------------
$shp = $s.Shapes.Item(9)
$shp.Left = 395.5443420410156
$shp.Top = 318.48095703125

# ---------------------------------------------------------------------
# Shape 10 (id 20) "TextBox 19" - "read write"
# ---------------------------------------------------------------------
$shp = $s.Shapes.Item(10)
$tr = $shp.TextFrame.TextRange
$tr.Font.Size = 12
$tr.Font.Name = "IBM Plex Sans"
$shp.Left = 435.5581970214844
$shp.Top = 296.41937255859375
$shp.Width = 69.82527923583984
$shp.Height = 21.810945510864258

# ---------------------------------------------------------------------
# Shape 11 (id 21) "TextBox 20" - two-paragraph caption
# ---------------------------------------------------------------------
$shp = $s.Shapes.Item(11)
$tr = $shp.TextFrame.TextRange
$para1 = $tr.Paragraphs(1)
$para1.Font.Size = 14
$para1.Font.Name = "IBM Plex Sans"
$para2 = $tr.Paragraphs(2)
$para2.Font.Size = 14
$para2.Font.Name = "IBM Plex Sans"
$shp.Left = 245.16458129882812
$shp.Top = 235.83346557617188
$shp.Width = 133.11859130859375
$shp.Height = 75.12654113769531
